# Apply the edits described by the commit diff across the three sheets:
# Summary, Assets, Liabilities.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Summary
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B3").Value = "Reem Al Awani"     # Name
$wsSummary.Range("B4").Value = 2797.99             # Monthly Income (AED)
$wsSummary.Range("B6").Value = 780155              # Total Assets (AED)
$wsSummary.Range("B7").Value = 30623               # Total Liabilities (AED)
$wsSummary.Range("B8").Value = 749532              # Net Worth (AED)
$wsSummary.Range("B9").Value = 25.48                # Asset/Liability Ratio

# ---------------------------------------------------------------------
# Sheet 2: Assets
# ---------------------------------------------------------------------
$wsAssets = $wb.Worksheets.Item("Assets")

# Insert two new rows above the existing "Liquid Assets" row (row 2),
# pushing the existing data (Liquid Assets row, TOTAL ASSETS row) down
# by two rows.
$wsAssets.Rows.Item(2).Insert()
$wsAssets.Rows.Item(2).Insert()

# New row 2: Vehicles / Luxury Car / 458952
$wsAssets.Range("A2").Value = "Vehicles"
$wsAssets.Range("B2").Value = "Luxury Car"
$wsAssets.Range("C2").Value = 458952

# New row 3: Vehicles / Luxury Car / 316219
$wsAssets.Range("A3").Value = "Vehicles"
$wsAssets.Range("B3").Value = "Luxury Car"
$wsAssets.Range("C3").Value = 316219

# Match the formatting of the surrounding data rows (style carried by the
# former row 2 / now row 4) for the two newly-inserted rows.
$wsAssets.Range("A4:C4").Copy()
$wsAssets.Range("A2:C3").PasteSpecial(-4122)

# Update values on the (now shifted) existing rows.
$wsAssets.Range("C4").Value = 4984      # Liquid Assets / Savings Account
$wsAssets.Range("C5").Value = 780155    # TOTAL ASSETS

# ---------------------------------------------------------------------
# Sheet 3: Liabilities
# ---------------------------------------------------------------------
$wsLiabilities = $wb.Worksheets.Item("Liabilities")

$wsLiabilities.Range("C2").Value = 30623   # Credit Card Balance amount
$wsLiabilities.Range("D2").Value = 1531    # Monthly Payment
$wsLiabilities.Range("C3").Value = 30623   # TOTAL LIABILITIES

Write-Host "Edits applied"
